# Repository.TestGitCommitChangesOutsideWebstudio/Main.xlsx
# "update file with jgit" - cell E8 ("Good Morning") was changed to "GIT UPDATE",
# and the active selection moved to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
